$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: D-column values look like plain numbers (or multi-dot strings) to Excel.
# Prefixing with a leading apostrophe forces a literal-text entry, exactly like a
# user typing an apostrophe before a numeric-looking value, keeping the General
# number format (only a quotePrefix flag is set) instead of silently becoming a number.

$ws.Range("D2").Value = "'42.222.07"
$ws.Range("E2").Value = "  -3.78%  "

$ws.Range("D3").Value = "'2.245.21"
$ws.Range("E3").Value = "  -4.56%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'232.54"
$ws.Range("E5").Value = "  -3.34%  "

$ws.Range("D6").Value = "'0.636"
$ws.Range("E6").Value = "  -5.57%  "

$ws.Range("D7").Value = "'70.74"
$ws.Range("E7").Value = "  -3.29%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "'0.562"
$ws.Range("E9").Value = "  -6.37%  "

$ws.Range("D10").Value = "'0.101"
$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("D11").Value = "'58.48"
$ws.Range("E11").Value = "  -1.22%  "

$ws.Range("D12").Value = "'36.07"
$ws.Range("E12").Value = "  +7.49%  "

$ws.Range("E13").Value = "  -2.38%  "

$ws.Range("D14").Value = "'6.88"
$ws.Range("E14").Value = "  -6.48%  "

$ws.Range("D15").Value = "'2.581.99"
$ws.Range("E15").Value = "  -4.38%  "

$ws.Range("D16").Value = "'15.14"
$ws.Range("E16").Value = "  -7.82%  "

$ws.Range("D17").Value = "'0.873"
$ws.Range("E17").Value = "  -3.93%  "

$ws.Range("D18").Value = "'2.247.25"
$ws.Range("E18").Value = "  -4.45%  "

$ws.Range("D19").Value = "'42.118.47"
$ws.Range("E19").Value = "  -3.77%  "

$ws.Range("D20").Value = "'0.0₃0997"
$ws.Range("E20").Value = "  -2.70%  "

$ws.Range("D21").Value = "'73.80"
$ws.Range("E21").Value = "  -4.66%  "

$ws.Range("D22").Value = "'6.23"
$ws.Range("E22").Value = "  -7.03%  "

$ws.Range("D23").Value = "'239.06"
$ws.Range("E23").Value = "  -6.96%  "

$ws.Range("D24").Value = "'1.96"
$ws.Range("E24").Value = "  +2.02%  "

$ws.Range("D26").Value = "'3.67"
$ws.Range("E26").Value = "  -1.62%  "

$ws.Range("D27").Value = "'2.35"
$ws.Range("E27").Value = "  -6.08%  "

$ws.Range("D28").Value = "'10.20"
$ws.Range("E28").Value = "  -3.48%  "

$ws.Range("E29").Value = "  -5.20%  "

$ws.Range("D30").Value = "'168.40"
$ws.Range("E30").Value = "  -5.29%  "

$ws.Range("D31").Value = "'20.78"
$ws.Range("E31").Value = "  -8.30%  "

$ws.Range("D32").Value = "'0.120"
$ws.Range("E32").Value = "  -6.71%  "

$ws.Range("D33").Value = "'0.127"
$ws.Range("E33").Value = "  -6.29%  "

$ws.Range("D34").Value = "'0.0728"
$ws.Range("E34").Value = "  -3.93%  "

$ws.Range("D35").Value = "'5.37"
$ws.Range("E35").Value = "  -1.05%  "

$ws.Range("D36").Value = "'4.81"
$ws.Range("E36").Value = "  -7.53%  "

$ws.Range("D37").Value = "'3.62"
$ws.Range("E37").Value = "  -5.41%  "

$ws.Range("E38").Value = "  +18.34%  "

$ws.Range("D39").Value = "'6.13"
$ws.Range("E39").Value = "  -3.75%  "

$ws.Range("D40").Value = "'2.25"
$ws.Range("E40").Value = "  -6.00%  "

$ws.Range("D41").Value = "'0.0267"
$ws.Range("E41").Value = "  -4.07%  "

$ws.Range("D42").Value = "'67.52"
$ws.Range("E42").Value = "  +1.33%  "

$ws.Range("D43").Value = "'4.94"
$ws.Range("E43").Value = "  -2.60%  "

$ws.Range("D44").Value = "'8.94"
$ws.Range("E44").Value = "  -2.01%  "

$ws.Range("D45").Value = "'0.101"
$ws.Range("E45").Value = "  -9.24%  "

$ws.Range("D48").Value = "'10.39"
$ws.Range("E48").Value = "  +9.13%  "

$ws.Range("D49").Value = "'4.43"
$ws.Range("E49").Value = "  +4.97%  "

$ws.Range("D50").Value = "'2.37"
$ws.Range("E50").Value = "  -5.20%  "

$ws.Range("E51").Value = "  -6.07%  "

# Row 46/47: Algorand and BinanceUSD swap positions with updated values
$ws.Range("B46").Value = "BinanceUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.18%  "

$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "'0.189"
$ws.Range("E47").Value = "  -6.24%  "
